$wb = $excel.ActiveWorkbook

# Row 43 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 699.7222
$ws.Range("I43").Value = 692.2
$ws.Range("J43").Value = 702.61536
$ws.Range("K43").Value = 692.2
$ws.Range("L43").Value = 702.61536
$ws.Range("M43").Value = -623.2
$ws.Range("N43").Value = -840.61536

# Row 86 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1763.8
$ws.Range("I86").Value = 1587.875
$ws.Range("J86").Value = 1964.8572
$ws.Range("K86").Value = 1587.875
$ws.Range("L86").Value = 1964.8572
$ws.Range("M86").Value = -464.875
$ws.Range("N86").Value = -4210.8572

# Row 89 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1763.8
$ws.Range("I89").Value = 1587.875
$ws.Range("J89").Value = 1964.8572
$ws.Range("K89").Value = 7939.375
$ws.Range("L89").Value = 9824.286
$ws.Range("M89").Value = -2323.375
$ws.Range("N89").Value = -21056.286

# Row 106 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 222228300
$ws.Range("I106").Value = 66673970
$ws.Range("K106").Value = 66673970
$ws.Range("M106").Value = -66673339

# Row 112 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 11905708
$ws.Range("J112").Value = 14286682
$ws.Range("L112").Value = 42860046
$ws.Range("N112").Value = -42862262

# Row 129 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 832.8570999999999
$ws.Range("I129").Value = 439
$ws.Range("J129").Value = 911.6286
$ws.Range("K129").Value = 1317
$ws.Range("L129").Value = 2734.8858
$ws.Range("M129").Value = 3683
$ws.Range("N129").Value = -12734.8858

# Row 131 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 3141.4285
$ws.Range("I131").Value = 1045
$ws.Range("J131").Value = 3980
$ws.Range("K131").Value = 3135
$ws.Range("L131").Value = 11940
$ws.Range("M131").Value = 1905
$ws.Range("N131").Value = -22020

# Row 133 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 26333.334
$ws.Range("J133").Value = 26333.334
$ws.Range("L133").Value = 26333.334
$ws.Range("N133").Value = -36453.334

# Row 135 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 3190.9744
$ws.Range("I135").Value = 2424.4443
$ws.Range("J135").Value = 4915.6665
$ws.Range("K135").Value = 21819.9987
$ws.Range("L135").Value = 44240.9985
$ws.Range("M135").Value = -19284.9987
$ws.Range("N135").Value = -49310.9985

# Row 32 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5155.532
$ws.Range("I32").Value = 4186.3267
$ws.Range("K32").Value = 4186.3267
$ws.Range("M32").Value = -3899.3267

# Row 45 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9954.77
$ws.Range("I45").Value = 14926.5
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 14926.5
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -14549.5
$ws.Range("N45").Value = -2754

# Row 132 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4578.375
$ws.Range("I132").Value = 1477.8158
$ws.Range("J132").Value = 16360.5
$ws.Range("K132").Value = 4433.4474
$ws.Range("L132").Value = 49081.5
$ws.Range("M132").Value = -1903.4474
$ws.Range("N132").Value = -54141.5

# Row 134 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# Row 86 on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2352.6924
$ws.Range("I86").Value = 2247.8
$ws.Range("J86").Value = 2702.3333
$ws.Range("K86").Value = 2247.8
$ws.Range("L86").Value = 2702.3333
$ws.Range("M86").Value = -1124.8
$ws.Range("N86").Value = -4948.3333

# Row 89 on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2352.6924
$ws.Range("I89").Value = 2247.8
$ws.Range("J89").Value = 2702.3333
$ws.Range("K89").Value = 11239
$ws.Range("L89").Value = 13511.6665
$ws.Range("M89").Value = -5623
$ws.Range("N89").Value = -24743.6665

# Row 105 on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2139.9443
$ws.Range("I105").Value = 1893.5
$ws.Range("J105").Value = 3002.5
$ws.Range("K105").Value = 1893.5
$ws.Range("L105").Value = 3002.5
$ws.Range("M105").Value = -146.5
$ws.Range("N105").Value = -6496.5

# Row 31 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6880.5
$ws.Range("I31").Value = 1154.8148
$ws.Range("K31").Value = 1154.8148
$ws.Range("M31").Value = -859.8148000000001

# Row 34 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6880.5
$ws.Range("I34").Value = 1154.8148
$ws.Range("K34").Value = 1154.8148
$ws.Range("M34").Value = -952.8148000000001

# Row 99 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 15627041
$ws.Range("I99").Value = 2332.2856
$ws.Range("K99").Value = 2332.2856
$ws.Range("M99").Value = -834.2856000000002

# Row 105 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1614.1666
$ws.Range("I105").Value = 1614.1666
$ws.Range("K105").Value = 1614.1666
$ws.Range("M105").Value = 132.8334

# Row 124 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H124").Value = 39800
$ws.Range("J124").Value = 39800
$ws.Range("L124").Value = 39800
$ws.Range("N124").Value = -44710

# Row 126 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 15627041
$ws.Range("I126").Value = 2332.2856
$ws.Range("K126").Value = 6996.8568
$ws.Range("M126").Value = -4526.8568

# Row 127 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# Row 132 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1674.0286
$ws.Range("I132").Value = 1266.0333
$ws.Range("J132").Value = 4122
$ws.Range("K132").Value = 3798.0999
$ws.Range("L132").Value = 12366
$ws.Range("M132").Value = -1268.0999
$ws.Range("N132").Value = -17426

# Row 104 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 6166
$ws.Range("I104").Value = 5998
$ws.Range("J104").Value = 6250
$ws.Range("K104").Value = 17994
$ws.Range("L104").Value = 18750
$ws.Range("M104").Value = -15373
$ws.Range("N104").Value = -23992

# Row 122 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 780.8
$ws.Range("J122").Value = 800
$ws.Range("L122").Value = 7200
$ws.Range("N122").Value = -12100

# Row 138 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 11933.077
$ws.Range("I138").Value = 14347.777
$ws.Range("K138").Value = 43043.331
$ws.Range("M138").Value = -37903.331

# Row 139 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 3932.0833
$ws.Range("I139").Value = 4928
$ws.Range("K139").Value = 14784
$ws.Range("M139").Value = -9644

# Row 102 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3647.7896
$ws.Range("I102").Value = 2129.6
$ws.Range("J102").Value = 5334.6665
$ws.Range("K102").Value = 2129.6
$ws.Range("L102").Value = 5334.6665
$ws.Range("M102").Value = -507.5999999999999
$ws.Range("N102").Value = -8578.666499999999

# Row 122 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4988284.5
$ws.Range("I122").Value = 5894154.5
$ws.Range("J122").Value = 5999.5
$ws.Range("K122").Value = 17682463.5
$ws.Range("L122").Value = 17998.5
$ws.Range("M122").Value = -17680013.5
$ws.Range("N122").Value = -22898.5

# Row 132 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2186.0159
$ws.Range("I132").Value = 1586.3684
$ws.Range("J132").Value = 3097.48
$ws.Range("K132").Value = 4759.1052
$ws.Range("L132").Value = 9292.440000000001
$ws.Range("M132").Value = -2229.1052
$ws.Range("N132").Value = -14352.44

# Row 26 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 14990
$ws.Range("I26").Value = 14990
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 14990
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -14695
$ws.Range("N26").ClearContents()

# Row 40 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 71432110
$ws.Range("I40").Value = 90912780
$ws.Range("K40").Value = 90912780
$ws.Range("M40").Value = -90912644

# Row 82 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1383225.2
$ws.Range("I82").Value = 2501200.5
$ws.Range("J82").Value = 265250
$ws.Range("K82").Value = 2501200.5
$ws.Range("L82").Value = 265250
$ws.Range("M82").Value = -2500839.5
$ws.Range("N82").Value = -265972

# Row 85 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1383225.2
$ws.Range("I85").Value = 2501200.5
$ws.Range("J85").Value = 265250
$ws.Range("K85").Value = 2501200.5
$ws.Range("L85").Value = 265250
$ws.Range("M85").Value = -2499952.5
$ws.Range("N85").Value = -267746

# Row 132 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9224326
$ws.Range("I132").Value = 14947874
$ws.Range("J132").Value = 3054.7222
$ws.Range("K132").Value = 44843622
$ws.Range("L132").Value = 9164.1666
$ws.Range("M132").Value = -44841092
$ws.Range("N132").Value = -14224.1666

# Row 136 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 11307.708
$ws.Range("I136").Value = 10093.733
$ws.Range("J136").Value = 13331
$ws.Range("K136").Value = 30281.199
$ws.Range("L136").Value = 39993
$ws.Range("M136").Value = -27731.199
$ws.Range("N136").Value = -45093

# Row 33 on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 40000
$ws.Range("J33").Value = 40000
$ws.Range("L33").Value = 40000
$ws.Range("N33").Value = -40500

# Row 36 on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H36").Value = 40000
$ws.Range("J36").Value = 40000
$ws.Range("L36").Value = 40000
$ws.Range("N36").Value = -40500

# Row 126 on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1345.1666
$ws.Range("I126").Value = 854.53845
$ws.Range("J126").Value = 2620.8
$ws.Range("K126").Value = 2563.61535
$ws.Range("L126").Value = 7862.400000000001
$ws.Range("M126").Value = -93.61535000000003
$ws.Range("N126").Value = -12802.4
